$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Source data refresh: the "LatestPeriod" column (B) for the APS-derived rows
# (employment / self-employment / unemployment / inactivity rate + volume
# rows, 2-9) moves on from "Jan 2024-Dec 2024" to the newer APS release
# window "Apr 2024 - Mar 2025".
foreach ($r in 2..9) {
    $ws.Cells.Item($r, 2).Value = "Apr 2024 - Mar 2025"
}

# Reflect the new scroll position / selection left after the edits (the
# sheet had been scrolled down to review row 10; now the top of the table
# is back in view with F4 selected).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F4").Select()
